$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '67.238.40'
$ws.Range("E2").Value = '  +0.77%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '3.472.88'

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '592.59'
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '178.18'
$ws.Range("E6").Value = '  +4.03%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '3.473.57'
$ws.Range("E8").Value = '  -0.40%  '

$ws.Range("E9").Value = '  -0.58%  '

$ws.Range("E10").Value = '  +5.12%  '

$ws.Range("E11").Value = '  -2.44%  '

$ws.Range("E12").Value = '  +0.33%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '4.072.43'
$ws.Range("E13").Value = '  -0.44%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '31.92'
$ws.Range("E14").Value = '  +11.14%  '

$ws.Range("E15").Value = '  +1.61%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '67.256.10'
$ws.Range("E16").Value = '  +0.67%  '

$ws.Range("E17").Value = '  -0.03%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '3.474.77'
$ws.Range("E18").Value = '  -0.26%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '6.23'
$ws.Range("E19").Value = '  -0.35%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '14.24'
$ws.Range("E20").Value = '  +1.62%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '387.93'
$ws.Range("E21").Value = '  -1.02%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '7.84'
$ws.Range("E22").Value = '  -0.82%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '72.99'
$ws.Range("E23").Value = '  +0.36%  '

$ws.Range("E24").Value = '  -0.23%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '5.71'
$ws.Range("E25").Value = '  +0.50%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '0.533'
$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '0.0000120'
$ws.Range("E27").Value = '  +0.85%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '10.33'
$ws.Range("E28").Value = '  +1.74%  '

$ws.Range("E29").Value = '  -3.13%  '

$ws.Range("E30").Value = '  -0.43%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '6.14'
$ws.Range("E31").Value = '  -0.58%  '

$ws.Range("E32").Value = '  -0.25%  '

$ws.Range("E33").Value = '  +0.17%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '23.46'
$ws.Range("E34").Value = '  -0.70%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '7.36'
$ws.Range("E35").Value = '  +0.80%  '

$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("E37").Value = '  -1.78%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '163.77'
$ws.Range("E38").Value = '  +0.60%  '

$ws.Range("E39").Value = '  -0.79%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '1.87'
$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '2.71'
$ws.Range("E41").Value = '  +6.41%  '

$ws.Range("E43").Value = '  -0.88%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '26.24'
$ws.Range("E44").Value = '  +0.88%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '2.819.61'
$ws.Range("E45").Value = '  +0.87%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.0719'
$ws.Range("E46").Value = '  -2.49%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '26.41'
$ws.Range("E47").Value = '  -2.50%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '41.50'
$ws.Range("E48").Value = '  -2.66%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.0298'
$ws.Range("E49").Value = '  -0.81%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '336.00'
$ws.Range("E50").Value = '  +0.01%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '1.05'
$ws.Range("E51").Value = '  -2.53%  '
